$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet: 展览 (Exhibition)
$ws1.Range("F2").Value = 1916
$ws1.Range("F3").Value = 517
$ws1.Range("F6").Value = 2807
$ws1.Range("F7").Value = 195
$ws1.Range("F8").Value = 100
$ws1.Range("F9").Value = 184
$ws1.Range("F10").Value = 1584
$ws1.Range("F11").Value = 563
$ws1.Range("F13").Value = 343
$ws1.Range("F17").Value = 13
$ws1.Range("F22").Value = 2
$ws1.Range("F24").Value = 239
$ws1.Range("F25").Value = 25
$ws1.Range("F27").Value = 1791
$ws1.Range("F28").Value = 42
$ws1.Range("F29").Value = 429
$ws1.Range("F30").Value = 94
$ws1.Range("F34").Value = 461

# Sheet: 全部类型 (All types)
$ws4.Range("F2").Value = 1916
$ws4.Range("F4").Value = 517
$ws4.Range("F7").Value = 2807
$ws4.Range("F8").Value = 195
$ws4.Range("F9").Value = 100
$ws4.Range("F10").Value = 184
$ws4.Range("F11").Value = 1584
$ws4.Range("F12").Value = 563
$ws4.Range("F14").Value = 343
$ws4.Range("F18").Value = 13
$ws4.Range("F23").Value = 2
$ws4.Range("F25").Value = 239
$ws4.Range("F26").Value = 25
$ws4.Range("F28").Value = 1791
$ws4.Range("F29").Value = 42
$ws4.Range("F30").Value = 429
$ws4.Range("F31").Value = 94
$ws4.Range("F35").Value = 461
